$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '46.689.85'
$ws.Range('E2').Value = '  +6.26%  '
$ws.Range('D3').Value = '2.298.94'
$ws.Range('E3').Value = '  +3.28%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').Value = '304.97'
$ws.Range('E5').Value = '  +1.93%  '
$ws.Range('D6').Value = '101.47'
$ws.Range('E6').Value = '  +12.23%  '
$ws.Range('D7').Value = '0.569'
$ws.Range('E7').Value = '  +2.75%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').Value = '0.524'
$ws.Range('E9').Value = '  +6.72%  '
$ws.Range('D10').Value = '36.89'
$ws.Range('E10').Value = '  +11.81%  '
$ws.Range('E11').Value = '  +2.35%  '
$ws.Range('D12').Value = '7.48'
$ws.Range('E12').Value = '  +7.83%  '
$ws.Range('E13').Value = '  -0.02%  '
$ws.Range('D14').Value = '2.650.71'
$ws.Range('E14').Value = '  +3.33%  '
$ws.Range('D15').Value = '2.300.69'
$ws.Range('E15').Value = '  +3.33%  '
$ws.Range('D16').Value = '13.97'
$ws.Range('E16').Value = '  +3.92%  '
$ws.Range('E17').Value = '  +5.54%  '
$ws.Range('D18').Value = '46.693.29'
$ws.Range('E18').Value = '  +6.63%  '
$ws.Range('D19').Value = '13.33'
$ws.Range('E19').Value = '  +17.82%  '
$ws.Range('D20').Value = '0.0₃0949'
$ws.Range('E20').Value = '  +5.08%  '
$ws.Range('D21').Value = '6.07'
$ws.Range('E21').Value = '  +2.01%  '
$ws.Range('D22').Value = '66.78'
$ws.Range('E22').Value = '  +3.34%  '
$ws.Range('D23').Value = '250.49'
$ws.Range('E23').Value = '  +6.04%  '
$ws.Range('D24').Value = '2.93'
$ws.Range('E24').Value = '  +4.32%  '
$ws.Range('E25').Value = '  +4.60%  '
$ws.Range('E26').Value = '  +0.11%  '
$ws.Range('D27').Value = '43.91'
$ws.Range('E27').Value = '  +15.03%  '
$ws.Range('E28').Value = '  +4.28%  '
$ws.Range('D29').Value = '9.93'
$ws.Range('E29').Value = '  +6.53%  '
$ws.Range('D30').Value = '20.16'
$ws.Range('E30').Value = '  +4.89%  '
$ws.Range('D31').Value = '5.73'
$ws.Range('E31').Value = '  +6.02%  '
$ws.Range('E32').Value = '  +11.57%  '
$ws.Range('D33').Value = '147.36'
$ws.Range('E33').Value = '  -2.96%  '
$ws.Range('D34').Value = '0.0799'
$ws.Range('E34').Value = '  +6.19%  '
$ws.Range('D35').Value = '3.22'
$ws.Range('E35').Value = '  +14.16%  '
$ws.Range('E36').Value = '  +11.65%  '
$ws.Range('E37').Value = '  +3.12%  '
$ws.Range('E38').Value = '  +5.22%  '
$ws.Range('D39').Value = '16.10'
$ws.Range('E39').Value = '  +22.05%  '
$ws.Range('D40').Value = '4.11'
$ws.Range('E40').Value = '  +12.73%  '
$ws.Range('D41').Value = '3.43'
$ws.Range('E41').Value = '  +7.26%  '
$ws.Range('D42').Value = '0.0304'
$ws.Range('E42').Value = '  +1.41%  '
$ws.Range('E43').Value = '  -0.03%  '
$ws.Range('D44').Value = '1.98'
$ws.Range('E44').Value = '  +11.79%  '
$ws.Range('D45').Value = '1.854.36'
$ws.Range('E45').Value = '  +0.98%  '
$ws.Range('D46').Value = '88.21'
$ws.Range('E46').Value = '  +20.45%  '
$ws.Range('D47').Value = '0.196'
$ws.Range('E47').Value = '  +8.45%  '
$ws.Range('D48').Value = '74.21'
$ws.Range('E48').Value = '  +9.38%  '
$ws.Range('E49').Value = '  +12.26%  '
$ws.Range('D50').Value = '96.60'
$ws.Range('E50').Value = '  +2.64%  '
$ws.Range('D51').Value = '54.14'
$ws.Range('E51').Value = '  +6.88%  '
